$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three changed data values in column C
$ws.Range("C2").Value = 58
$ws.Range("C3").Value = 65
$ws.Range("C17").Value = 55

# Re-apply alignment to the whole data column so Excel collapses the
# duplicate/redundant style into the shared "centered" style used
# elsewhere on the sheet (matches cellXfs count 4 -> 3 in the diff).
$dataRange = $ws.Range("C2:C49")
$dataRange.HorizontalAlignment = -4108  # xlCenter
$dataRange.VerticalAlignment = -4108    # xlCenter

# Update the active selection to match the new cursor position
$ws.Range("D3").Select()
